# feat: agregar script para actualizar DocEntries desde Excel
#
# Replaces the 11 numeric DocEntry sample values (A2:A12) with a new list
# of 46 DocEntry values (A2:A47), written as text so leading structure is
# preserved exactly as exported from the source system, then updates the
# sheet view to match (scrolled to row 11, A2:A47 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docEntries = @(
    "123729", "154974", "65891",  "114164", "152122",
    "151116", "102190", "114446", "117997", "147340",
    "98679",  "103452", "93505",  "120498", "112719",
    "93653",  "119106", "147336", "108080", "101189",
    "108308", "101298", "147322", "101284", "101287",
    "93654",  "101288", "101274", "101276", "101279",
    "152669", "143647", "154190", "101745", "109324",
    "104309", "140120", "147228", "156677", "162039",
    "101775", "142221", "149041", "147932", "157593",
    "162505"
)

# Clear out the old sample rows (A2:A12) before writing the new, longer list.
$ws.Range("A2:A12").ClearContents()

$lastRow = 1 + $docEntries.Count
$row = 2
foreach ($entry in $docEntries) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $entry
    $row++
}

# Keep the header row's layout untouched; just resize the view/selection to
# the new extent of the data.
$ws.Range("A2:A$lastRow").NumberFormat = "@"

$excel.ActiveWindow.ScrollRow = 11
$ws.Range("A2:A$lastRow").Select()
